$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff performs a 3-way cyclic rotation of the "Fecha" (D), "Volumen" (M),
# "Precio minimo" (N), "Precio maximo" (O), "Precio promedio ponderado" (P),
# "Origen" (R) and "Precio $/Kg" (S) values across rows 2, 4 and 5:
#   new row2 = old row5, new row4 = old row2, new row5 = old row4

$cols = @("D", "M", "N", "O", "P", "R", "S")

# Capture the original values before any writes.
$orig2 = @{}
$orig4 = @{}
$orig5 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("${col}2").Value2
    $orig4[$col] = $ws.Range("${col}4").Value2
    $orig5[$col] = $ws.Range("${col}5").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $orig5[$col]
    $ws.Range("${col}4").Value2 = $orig2[$col]
    $ws.Range("${col}5").Value2 = $orig4[$col]
}
